# Append one new row (row 74) of sensor data to each of the four worksheets,
# mirroring the prior row (row 73) but with the timestamp advanced by one hour.

$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{ Name = "ROW35-FE-LIFTER";  Time = "2025-03-07 09:42:06"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; E = "0x d";  G = "568631262647113770877196"; I = 13 },
    @{ Name = "ROW35-MID-LIFTER"; Time = "2025-03-07 09:29:35"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; E = "0x e";  G = "568631262647113770942732"; I = 14 },
    @{ Name = "ROW02-FE-LIFTER";  Time = "2025-03-07 09:51:45"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; E = "0xff";  G = "568631262647113769959692"; I = 255 },
    @{ Name = "ROW02-MID-LIFTER"; Time = "2025-03-07 09:41:15"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; E = "0x 3";  G = "568631262647113769959692"; I = 3 }
)

foreach ($sd in $sheetsData) {
    $ws = $wb.Worksheets.Item($sd.Name)
    $row = 74

    $ws.Cells.Item($row, 1).Value = $sd.Time
    $ws.Cells.Item($row, 2).Value = "0x01,0x90 "
    $ws.Cells.Item($row, 3).Value = $sd.C
    $ws.Cells.Item($row, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($row, 5).Value = $sd.E
    $ws.Cells.Item($row, 6).Value = 400
    # Column G holds a long digit-string that must stay text (it would lose
    # precision as a double), so force text formatting before assigning it.
    $gCell = $ws.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $sd.G
    $ws.Cells.Item($row, 8).Value = 400
    $ws.Cells.Item($row, 9).Value = $sd.I
}
